$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric values for rows 2-6 (per diff)
# Row 2
$ws.Range("D2").Value2 = 1895
$ws.Range("E2").Value2 = -3
$ws.Range("F2").Value2 = -3
$ws.Range("G2").Value2 = 11
$ws.Range("H2").Value2 = 11
$ws.Range("I2").Value2 = 11
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 1983
$ws.Range("L2").Value2 = 776
$ws.Range("M2").Value2 = 1206
$ws.Range("N2").Value2 = 1208
$ws.Range("O2").Value2 = -2
$ws.Range("P2").Value2 = 143
$ws.Range("Q2").Value2 = 93
$ws.Range("R2").Value2 = 255
$ws.Range("S2").Value2 = -216
$ws.Range("T2").Value2 = 13
$ws.Range("U2").Value2 = 80
$ws.Range("V2").Value2 = 110
$ws.Range("W2").Value2 = -0.18
$ws.Range("X2").Value2 = 0.58
$ws.Range("Y2").Value2 = 0.95
$ws.Range("Z2").Value2 = 0.51
$ws.Range("AA2").Value2 = 64.36
$ws.Range("AB2").Value2 = 730.03
$ws.Range("AC2").Value2 = 32
$ws.Range("AD2").Value2 = 27.19
$ws.Range("AE2").Value2 = 3624
$ws.Range("AF2").Value2 = 0.24
$ws.Range("AG2").Value2 = 0
$ws.Range("AH2").Value2 = 0
$ws.Range("AI2").Value2 = 0
$ws.Range("AJ2").Value2 = 36274078

# Row 3
$ws.Range("D3").Value2 = 1324
$ws.Range("E3").Value2 = -23
$ws.Range("F3").Value2 = -14
$ws.Range("G3").Value2 = -23
$ws.Range("H3").Value2 = -84
$ws.Range("I3").Value2 = -83
$ws.Range("J3").Value2 = -1
$ws.Range("K3").Value2 = 2172
$ws.Range("L3").Value2 = 977
$ws.Range("M3").Value2 = 1196
$ws.Range("N3").Value2 = 1198
$ws.Range("O3").Value2 = -2
$ws.Range("P3").Value2 = 234
$ws.Range("Q3").Value2 = 18
$ws.Range("R3").Value2 = -243
$ws.Range("S3").Value2 = 63
$ws.Range("T3").Value2 = 2
$ws.Range("U3").Value2 = 16
$ws.Range("V3").Value2 = 77
$ws.Range("W3").Value2 = -1.72
$ws.Range("X3").Value2 = -6.31
$ws.Range("Y3").Value2 = -6.87
$ws.Range("Z3").Value2 = -4.02
$ws.Range("AA3").Value2 = 81.67
$ws.Range("AB3").Value2 = 405.56
$ws.Range("AC3").Value2 = -201
$ws.Range("AD3").Value2 = -4.96
$ws.Range("AE3").Value2 = 2732
$ws.Range("AF3").Value2 = 0.36
$ws.Range("AG3").Value2 = 0
$ws.Range("AH3").Value2 = 0
$ws.Range("AI3").Value2 = 0
$ws.Range("AJ3").Value2 = 46803136

# Row 4
$ws.Range("D4").Value2 = 2409
$ws.Range("E4").Value2 = 18
$ws.Range("F4").Value2 = -22
$ws.Range("G4").Value2 = -17
$ws.Range("H4").Value2 = -113
$ws.Range("I4").Value2 = -129
$ws.Range("J4").Value2 = 16
$ws.Range("K4").Value2 = 4246
$ws.Range("L4").Value2 = 1813
$ws.Range("M4").Value2 = 2433
$ws.Range("N4").Value2 = 1177
$ws.Range("O4").Value2 = 1256
$ws.Range("P4").Value2 = 234
$ws.Range("Q4").Value2 = 127
$ws.Range("R4").Value2 = -9
$ws.Range("S4").Value2 = 120
$ws.Range("T4").Value2 = 207
$ws.Range("U4").Value2 = -81
$ws.Range("V4").Value2 = 1161
$ws.Range("W4").Value2 = 0.76
$ws.Range("X4").Value2 = -4.69
$ws.Range("Y4").Value2 = -10.82
$ws.Range("Z4").Value2 = -3.52
$ws.Range("AA4").Value2 = 74.51
$ws.Range("AB4").Value2 = 347.01
$ws.Range("AC4").Value2 = -275
$ws.Range("AD4").Value2 = -3.39
$ws.Range("AE4").Value2 = 2684
$ws.Range("AF4").Value2 = 0.35
$ws.Range("AG4").Value2 = 0
$ws.Range("AH4").Value2 = 0
$ws.Range("AI4").Value2 = 0
$ws.Range("AJ4").Value2 = 46803136

# Row 5
$ws.Range("D5").Value2 = 2464
$ws.Range("E5").Value2 = 62
$ws.Range("F5").Value2 = 62
$ws.Range("G5").Value2 = 41
$ws.Range("H5").Value2 = 50
$ws.Range("I5").Value2 = 44
$ws.Range("J5").Value2 = 6
$ws.Range("K5").Value2 = 4048
$ws.Range("L5").Value2 = 1423
$ws.Range("M5").Value2 = 2625
$ws.Range("N5").Value2 = 1194
$ws.Range("O5").Value2 = 1431
$ws.Range("P5").Value2 = 234
$ws.Range("Q5").Value2 = 158
$ws.Range("R5").Value2 = -87
$ws.Range("S5").Value2 = -209
$ws.Range("T5").Value2 = 66
$ws.Range("U5").Value2 = 92
$ws.Range("V5").Value2 = 792
$ws.Range("W5").Value2 = 2.53
$ws.Range("X5").Value2 = 2.02
$ws.Range("Y5").Value2 = 3.72
$ws.Range("Z5").Value2 = 1.2
$ws.Range("AA5").Value2 = 54.18
$ws.Range("AB5").Value2 = 355.25
$ws.Range("AC5").Value2 = 94
$ws.Range("AD5").Value2 = 8.55
$ws.Range("AE5").Value2 = 2723
$ws.Range("AF5").Value2 = 0.3
$ws.Range("AG5").Value2 = 20
$ws.Range("AH5").Value2 = 2.48
$ws.Range("AI5").Value2 = 19.86
$ws.Range("AJ5").Value2 = 46803136

# Row 6
$ws.Range("D6").Value2 = 2342
$ws.Range("E6").Value2 = 32
$ws.Range("F6").Value2 = 32
$ws.Range("G6").Value2 = 10
$ws.Range("H6").Value2 = 11
$ws.Range("I6").Value2 = -2
$ws.Range("K6").Value2 = 3908
$ws.Range("L6").Value2 = 1267
$ws.Range("M6").Value2 = 2641
$ws.Range("N6").Value2 = 1182
$ws.Range("P6").Value2 = 234
$ws.Range("Q6").Value2 = 247
$ws.Range("R6").Value2 = -113
$ws.Range("S6").Value2 = -111
$ws.Range("T6").Value2 = 89
$ws.Range("U6").Value2 = 158
$ws.Range("V6").Value2 = 681
$ws.Range("W6").Value2 = 1.38
$ws.Range("X6").Value2 = 0.45
$ws.Range("Y6").Value2 = -0.2
$ws.Range("Z6").Value2 = 0.27
$ws.Range("AA6").Value2 = 47.98
$ws.Range("AB6").Value2 = 347.75
$ws.Range("AC6").Value2 = -5
$ws.Range("AD6").Value2 = -162.04
$ws.Range("AE6").Value2 = 2695
$ws.Range("AF6").Value2 = 0.3
$ws.Range("AG6").Value2 = 15
$ws.Range("AH6").Value2 = 1.85
$ws.Range("AI6").Value2 = -281.17
$ws.Range("AJ6").Value2 = 46803136

# Rows 7-9: clear all data columns (D:AJ), keep only A, B, C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
